$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pokemon")

# --- Shift the existing "exp" column (H) one column to the right (I), ---
# --- keeping it intact, to make room for a new "vitesse" column in H.  ---
$ws.Range("I1").Value2 = $ws.Range("H1").Value2
$ws.Range("I2").Value2 = $ws.Range("H2").Value2
$ws.Range("I3").Value2 = $ws.Range("H3").Value2
$ws.Range("I4").Value2 = $ws.Range("H4").Value2

# --- New "vitesse" (speed) column in H with its values. ---
$ws.Range("H1").Value2 = "vitesse"
$ws.Range("H2").Value2 = 10
$ws.Range("H3").Value2 = 50
$ws.Range("H4").Value2 = 10

# --- Swap the "rarete" values between Pikachu (row 2) and Mewtwo (row 3). ---
# Use Copy/Paste (not Value=) so the text values aren't re-parsed as numbers.
$ws.Range("Z1").Value2 = $null
$ws.Range("Z2").Value2 = $null
$ws.Range("D2").Copy($ws.Range("Z1"))
$ws.Range("D3").Copy($ws.Range("Z2"))
$ws.Range("Z2").Copy($ws.Range("D2"))
$ws.Range("Z1").Copy($ws.Range("D3"))
$ws.Range("Z1:Z2").Clear()

# --- Make "pokemon" the active sheet / tab, with the new selection. ---
$ws.Activate()
$ws.Range("I9").Select()
